# minor clean up + initial visualizations
# Remove the "635.txt" rows (rows 67-71) and renumber the ID column (A)
# for the remaining data rows so it stays sequential.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 rows that hold the "635.txt" entries.
$ws.Rows("67:71").Delete()

# The ID column (A) is a plain sequential counter (1, 2, 3, ...).
# After the deletion, renumber it so it stays sequential through the
# last remaining data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells($r, 1).Value = $r - 1
}
